$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("0039-68"), shifting remaining rows up.
$ws.Rows(2).Delete()

# Now rows are: 1 Placa, 2 HHJP-31, 3 THO-09171, 4 IFJ-14061, 5 KBN-08251, 6 XYZ-0001
# Delete rows 3 through 6 (THO-09171, IFJ-14061, KBN-08251, XYZ-0001)
$ws.Range("A3:A6").EntireRow.Delete()

# Match the resulting selection from the source workbook.
[void]$ws.Range("C6").Select()
